# "Add basic box model"
#
# 1. Delete the standalone "New Game" textbox (TextBox 8 / id 9).
# 2. Re-purpose the big empty header textbox (TextBox 28 / id 29): move +
#    resize it into a single title-bar line and give it the combined
#    "Rock, Paper and Scissors ... New Game" caption.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) remove the now-redundant "New Game" box -----------------------
$newGameBox = $s.Shapes.Item("TextBox 8")
$newGameBox.Delete()

# --- 2) reflow the title textbox ---------------------------------------
$titleBox = $s.Shapes.Item("TextBox 28")

# EMU -> point conversion (1 pt = 1/72in = 12700 EMU). A tiny epsilon is
# added before the conversion so that, after PowerPoint's COM layer
# round-trips the value through a Single (float32), truncating back to
# EMU lands exactly on the target instead of one EMU short.
$emuToPt = 914400 / 72

$titleBox.Left   = (1637952 / $emuToPt) + 0.00005
$titleBox.Top    = (1110815 / $emuToPt) + 0.00005
$titleBox.Width  = (8623883 / $emuToPt) + 0.00005
$titleBox.Height = (369332  / $emuToPt) + 0.00005

$titleBox.TextFrame.TextRange.Text = "Rock, Paper and Scissors                `t`t`t`t`tNew Game"
